# Update MSME definitions: replace literal "<br/>" markers with real
# line breaks inside the relevant shared-string cells on the Summary sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

$ws.Range("B19").Value = "<5 Indstr. & trade`n<4 Services"
$ws.Range("D19").Value = "<610,000 Agriculture, `n<1,800,000 Industry, `n<2,400,000 Trade, `n<590,000 Services, `n<760,000 Construction"

$ws.Range("B20").Value = "<24 Indstr. `n<23 Trade`n<17 Services"
$ws.Range("D20").Value = "<4,100,000 Agriculture, `n<10,300,000 Industry, `n<14,000,000 Trade, `n<4,300,000 Services, `n<4,800,000 Construction"

$ws.Range("B21").Value = "<96 Indstr. `n<67 Trade`n<66 Services"
$ws.Range("D21").Value = "<24,100,000 Agriculture, `n<82,200,000 Industry, `n<111,900,000 Trade, `n<28,300,000 Services, `n<37,700,000 Construction"

$ws.Range("B22").Value = ">=96 Indstr.`n >=67 Trade`n>=66 Services"
$ws.Range("D22").Value = ">=24,100,000 Agriculture, `n>=82,200,000 Industry, `n>=111,900,000 Trade, `n>=28,300,000 Services, `n>=37,700,000 Construction"
